# Auto-generated edit script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix proveedor/razon-social names: stray commas used as separators were
#    meant to be periods (scraper formatting fix). Whole-cell, case-sensitive
#    replace so we only touch the exact strings below.
$ws.Cells.Replace('PALACIOS, GERMAN ALEJANDRO', 'PALACIOS. GERMAN ALEJANDRO', 1, 1, $true) | Out-Null
$ws.Cells.Replace('MARSICO GUILLERMO MIGUEL, MARSICO JUAN EDUARDO', 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO', 1, 1, $true) | Out-Null
$ws.Cells.Replace('RONDONI, LUCIANO MATIAS', 'RONDONI. LUCIANO MATIAS', 1, 1, $true) | Out-Null
$ws.Cells.Replace('SCHAB DARIO, PEROTTI XAVIER, BENINCA MATIAS S.H.', 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH', 1, 1, $true) | Out-Null
$ws.Cells.Replace('OLIVERA, FLORENCIO', 'OLIVERA. FLORENCIO', 1, 1, $true) | Out-Null
$ws.Cells.Replace('ALBIZZATTI, PABLO MARTIN Y FULINI, SERGIO RUBEN', 'ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN', 1, 1, $true) | Out-Null
$ws.Cells.Replace('BALDONI  ESTEBAN,BALDONI SEBASTIAN Y BALDONI LUCIANO S.H.', 'BALDONI  ESTEBAN.BALDONI SEBASTIAN Y BALDONI LUCIANO SH', 1, 1, $true) | Out-Null
$ws.Cells.Replace('DENING BLANCO, CRISTIAN DAVID', 'DENING BLANCO. CRISTIAN DAVID', 1, 1, $true) | Out-Null
$ws.Cells.Replace('PARRAVICINI VIRGINIA VANINA, VIRGINIA VANINA', 'PARRAVICINI VIRGINIA VANINA. VIRGINIA VANINA', 1, 1, $true) | Out-Null

# 2) Fix "Importe" amounts: values were stored as localized formatted text
#    (thousands separator ".", decimal comma ",") and must become plain
#    decimal-point numeric-looking text (no thousands separator, "." decimal).
#    Force the column to Text format first so Replace does not coerce the
#    result into a real number, then restore the original (default) style.
$importeRange = $ws.Range("H2:H280")
$importeRange.NumberFormat = "@"
$ws.Cells.Replace('30.000,00', '30000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('34.800,00', '34800.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('26.500,00', '26500.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('399.000,00', '399000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('712.569,00', '712569.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('1.103.333,33', '1103333.33', 1, 1, $true) | Out-Null
$ws.Cells.Replace('962.000,00', '962000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('980.000,00', '980000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('5.929,00', '5929.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('53.475,00', '53475.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('8.970,00', '8970.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('4.760,00', '4760.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('110.010,00', '110010.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('50.522,70', '50522.70', 1, 1, $true) | Out-Null
$ws.Cells.Replace('876.000,00', '876000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('1.087.999,50', '1087999.50', 1, 1, $true) | Out-Null
$ws.Cells.Replace('37.773,44', '37773.44', 1, 1, $true) | Out-Null
$ws.Cells.Replace('40.800,00', '40800.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('53.680,00', '53680.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('21.449,00', '21449.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('1.104.747,64', '1104747.64', 1, 1, $true) | Out-Null
$ws.Cells.Replace('510.766,43', '510766.43', 1, 1, $true) | Out-Null
$ws.Cells.Replace('248.569,60', '248569.60', 1, 1, $true) | Out-Null
$ws.Cells.Replace('776.138,18', '776138.18', 1, 1, $true) | Out-Null
$ws.Cells.Replace('19.181,25', '19181.25', 1, 1, $true) | Out-Null
$ws.Cells.Replace('3.510,00', '3510.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('1.569.431,36', '1569431.36', 1, 1, $true) | Out-Null
$ws.Cells.Replace('8.960,00', '8960.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('39.525,00', '39525.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('1.001.770,00', '1001770.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('70.592,60', '70592.60', 1, 1, $true) | Out-Null
$ws.Cells.Replace('62.149,39', '62149.39', 1, 1, $true) | Out-Null
$ws.Cells.Replace('115.874,86', '115874.86', 1, 1, $true) | Out-Null
$ws.Cells.Replace('6.300,00', '6300.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('25.000,00', '25000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('24.000,00', '24000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('3.150,00', '3150.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('11.000,00', '11000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('15.800,00', '15800.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('1.977,60', '1977.60', 1, 1, $true) | Out-Null
$ws.Cells.Replace('2.600,00', '2600.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('259,20', '259.20', 1, 1, $true) | Out-Null
$ws.Cells.Replace('138.436,29', '138436.29', 1, 1, $true) | Out-Null
$ws.Cells.Replace('21.885,00', '21885.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('600,00', '600.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('2.800,00', '2800.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('963,63', '963.63', 1, 1, $true) | Out-Null
$ws.Cells.Replace('83.220,00', '83220.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('1.600,00', '1600.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('23.698,72', '23698.72', 1, 1, $true) | Out-Null
$ws.Cells.Replace('6.850,12', '6850.12', 1, 1, $true) | Out-Null
$ws.Cells.Replace('13.245,00', '13245.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('22.000,00', '22000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('170.000,00', '170000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('195.007,36', '195007.36', 1, 1, $true) | Out-Null
$ws.Cells.Replace('80,00', '80.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('207.029,79', '207029.79', 1, 1, $true) | Out-Null
$ws.Cells.Replace('304.459,94', '304459.94', 1, 1, $true) | Out-Null
$ws.Cells.Replace('43.230,34', '43230.34', 1, 1, $true) | Out-Null
$ws.Cells.Replace('1.185,98', '1185.98', 1, 1, $true) | Out-Null
$ws.Cells.Replace('12.552,14', '12552.14', 1, 1, $true) | Out-Null
$ws.Cells.Replace('64.484,04', '64484.04', 1, 1, $true) | Out-Null
$ws.Cells.Replace('18.255,00', '18255.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('362,19', '362.19', 1, 1, $true) | Out-Null
$ws.Cells.Replace('231.867,45', '231867.45', 1, 1, $true) | Out-Null
$ws.Cells.Replace('835,36', '835.36', 1, 1, $true) | Out-Null
$ws.Cells.Replace('9.233,00', '9233.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('2.030,00', '2030.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('880.168,80', '880168.80', 1, 1, $true) | Out-Null
$ws.Cells.Replace('612,48', '612.48', 1, 1, $true) | Out-Null
$ws.Cells.Replace('13.387,80', '13387.80', 1, 1, $true) | Out-Null
$ws.Cells.Replace('3.600,00', '3600.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('35.676,64', '35676.64', 1, 1, $true) | Out-Null
$ws.Cells.Replace('760,00', '760.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('400,00', '400.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('6.482,00', '6482.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('7.200,00', '7200.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('28.379,11', '28379.11', 1, 1, $true) | Out-Null
$ws.Cells.Replace('25.825,00', '25825.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('6.720,00', '6720.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('8.027,70', '8027.70', 1, 1, $true) | Out-Null
$ws.Cells.Replace('6.122,60', '6122.60', 1, 1, $true) | Out-Null
$ws.Cells.Replace('61.700,00', '61700.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('16.985,50', '16985.50', 1, 1, $true) | Out-Null
$ws.Cells.Replace('800,00', '800.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('265.900,00', '265900.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('42.890,00', '42890.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('1.400,00', '1400.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('505.824,00', '505824.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('13.099,80', '13099.80', 1, 1, $true) | Out-Null
$ws.Cells.Replace('3.000,00', '3000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('24.835,00', '24835.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('27.830,00', '27830.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('140,40', '140.40', 1, 1, $true) | Out-Null
$ws.Cells.Replace('1.960,00', '1960.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('104.950,00', '104950.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('43.150,00', '43150.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('36.720,00', '36720.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('178.525,00', '178525.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('602.642,00', '602642.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('23.988,00', '23988.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('12.520,00', '12520.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('7,00', '7.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('5.541,27', '5541.27', 1, 1, $true) | Out-Null
$ws.Cells.Replace('28.140,00', '28140.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('3.500,00', '3500.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('28.700,00', '28700.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('19.690,00', '19690.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('62.940,00', '62940.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('276.620,00', '276620.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('4.985,00', '4985.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('21.000,00', '21000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('300,00', '300.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('7.800,00', '7800.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('905,00', '905.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('5.679,50', '5679.50', 1, 1, $true) | Out-Null
$ws.Cells.Replace('6.670,00', '6670.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('322,00', '322.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('9.800,00', '9800.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('12.681,00', '12681.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('4.520,00', '4520.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('175,00', '175.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('10.355,00', '10355.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('6.075,00', '6075.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('670,00', '670.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('172.137,39', '172137.39', 1, 1, $true) | Out-Null
$ws.Cells.Replace('955,00', '955.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('51.085,00', '51085.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('2.087,00', '2087.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('71.673,00', '71673.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('21.800,57', '21800.57', 1, 1, $true) | Out-Null
$ws.Cells.Replace('33.470,00', '33470.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('2.722,50', '2722.50', 1, 1, $true) | Out-Null
$ws.Cells.Replace('38.390,00', '38390.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('32.053,00', '32053.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('8.690,00', '8690.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('390,90', '390.90', 1, 1, $true) | Out-Null
$ws.Cells.Replace('16.940,00', '16940.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('28.900,00', '28900.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('7.600,00', '7600.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('30,00', '30.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('274,50', '274.50', 1, 1, $true) | Out-Null
$ws.Cells.Replace('4.400,00', '4400.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('18.200,00', '18200.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('9.000,00', '9000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('2.500,00', '2500.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('85.385,00', '85385.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('80.000,00', '80000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('14.000,00', '14000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('10.140,00', '10140.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('5.000,00', '5000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('3.931,37', '3931.37', 1, 1, $true) | Out-Null
$ws.Cells.Replace('103.378,58', '103378.58', 1, 1, $true) | Out-Null
$ws.Cells.Replace('749,24', '749.24', 1, 1, $true) | Out-Null
$ws.Cells.Replace('13.233,79', '13233.79', 1, 1, $true) | Out-Null
$ws.Cells.Replace('32.900,00', '32900.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('6.720,35', '6720.35', 1, 1, $true) | Out-Null
$ws.Cells.Replace('4.040,00', '4040.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('200,00', '200.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('16.500,00', '16500.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('20.000,00', '20000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('7.000,00', '7000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('10.000,00', '10000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('12.000,00', '12000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('12.500,00', '12500.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('10.350,00', '10350.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('10.500,00', '10500.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('15.000,00', '15000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('6.000,00', '6000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('45.000,02', '45000.02', 1, 1, $true) | Out-Null
$ws.Cells.Replace('13.800,00', '13800.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('4.000,00', '4000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('350.000,00', '350000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('6.480,00', '6480.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('4.950,00', '4950.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('8.000,00', '8000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('61.690,00', '61690.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('89.100,00', '89100.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('41.980,00', '41980.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('4.800,00', '4800.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('893,80', '893.80', 1, 1, $true) | Out-Null
$ws.Cells.Replace('1.338,67', '1338.67', 1, 1, $true) | Out-Null
$ws.Cells.Replace('42.842,48', '42842.48', 1, 1, $true) | Out-Null
$ws.Cells.Replace('2.978,80', '2978.80', 1, 1, $true) | Out-Null
$ws.Cells.Replace('14.374,80', '14374.80', 1, 1, $true) | Out-Null
$ws.Cells.Replace('18.637,29', '18637.29', 1, 1, $true) | Out-Null
$ws.Cells.Replace('17.018,00', '17018.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('15.352,69', '15352.69', 1, 1, $true) | Out-Null
$ws.Cells.Replace('61.100,00', '61100.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('9.890,00', '9890.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('84.950,00', '84950.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('20.004,00', '20004.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('5.290,00', '5290.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('24.133,00', '24133.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('24.500,00', '24500.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('19.600,00', '19600.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('12.780,00', '12780.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('3.550,00', '3550.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('87.024,05', '87024.05', 1, 1, $true) | Out-Null
$ws.Cells.Replace('40.000,00', '40000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('55.000,00', '55000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('15.730,00', '15730.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('267.000,00', '267000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('13.672,60', '13672.60', 1, 1, $true) | Out-Null
$ws.Cells.Replace('1.265.999,90', '1265999.90', 1, 1, $true) | Out-Null
$ws.Cells.Replace('476.000,00', '476000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('24.300,00', '24300.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('26.369.689,32', '26369689.32', 1, 1, $true) | Out-Null
$ws.Cells.Replace('445.720,00', '445720.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('387.020,00', '387020.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('409.350,00', '409350.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('355.000,00', '355000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('418.000,00', '418000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('678.000,00', '678000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('926.200,00', '926200.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('1.053.000,00', '1053000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('529.181,00', '529181.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('710.000,00', '710000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('747.280,00', '747280.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('701.100,00', '701100.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('1.043.000,00', '1043000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('725.320,00', '725320.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('1.041.860,00', '1041860.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('388.503,33', '388503.33', 1, 1, $true) | Out-Null
$ws.Cells.Replace('762.300,00', '762300.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('397.225,00', '397225.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('1.961.645,50', '1961645.50', 1, 1, $true) | Out-Null
$ws.Cells.Replace('1.508.684,23', '1508684.23', 1, 1, $true) | Out-Null
$ws.Cells.Replace('5.515.828,54', '5515828.54', 1, 1, $true) | Out-Null
$ws.Cells.Replace('20.757.741,28', '20757741.28', 1, 1, $true) | Out-Null
$ws.Cells.Replace('288.000,00', '288000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('779.000,00', '779000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('429.947,82', '429947.82', 1, 1, $true) | Out-Null
$ws.Cells.Replace('150.000,00', '150000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('161.000,00', '161000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('75.000,00', '75000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('21.900,00', '21900.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('180.207,50', '180207.50', 1, 1, $true) | Out-Null
$ws.Cells.Replace('147.468,75', '147468.75', 1, 1, $true) | Out-Null
$ws.Cells.Replace('7.970,00', '7970.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('19.000,00', '19000.00', 1, 1, $true) | Out-Null
$ws.Cells.Replace('26.540,00', '26540.00', 1, 1, $true) | Out-Null
$importeRange.Style = "Normal"

